# Remove fish from title slide
#
# The title slide (slide 1) has a small fish picture (Shape 38) layered
# over the subtitle, along with a fade-in click animation that targets it.
# Removing the shape leaves a dangling animation behind, so the animation
# effect is removed first and then the picture shape itself.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Drop the click animation effect(s) that target the fish picture before
# removing the shape, so we don't leave an orphaned <p:timing> entry.
$seq = $s.TimeLine.MainSequence
for ($i = $seq.Count; $i -ge 1; $i--) {
    $effect = $seq.Item($i)
    if ($effect.Shape.Name -eq "Shape 38") {
        $effect.Delete()
    }
}

# Now remove the fish picture shape itself.
foreach ($shp in $s.Shapes) {
    if ($shp.Name -eq "Shape 38") {
        $shp.Delete()
    }
}
